# Added Data to Investigation with ARCitect
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Investigation Identifier / Title / Description
$ws.Range("B7").Value = "M4.4_UC6_ARC"
$ws.Range("B8").Value = "ARC for UC6"
$ws.Range("B9").Value = "In this repository, we will create an annotated research context ARC with assays, studies, workflows and runs of these workflows for Use Case 6 of FAIRagro. Use Case 6 deals with ETL functions for semiautomated data integration into crop simulation modelling.
The ARC contians a modularised CWL workflow that can be executed locally with cwltool. To facilitate writing CWL, a minimal templating system that helps to generate CWL will be used or developed. Once the finalised ARC has undergone a test to verify that the workflow is running as intended, it will also be made accessible via the PLANTdataHUB. The ARC will serve as a test case for the further development of a Scientific Workflow Infrastructure (SciWIn)."

# INVESTIGATION CONTACTS
# Investigation Person Last Name (row 21)
$ws.Range("B21").Value = "Krumsieck"
$ws.Range("C21").Value = "Leidel"
$ws.Range("D21").Value = "König"
$ws.Range("E21").Value = "von Waldow"

# Investigation Person First Name (row 22)
$ws.Range("B22").Value = "Jens"
$ws.Range("C22").Value = "Antonia"
$ws.Range("D22").Value = "Patrick"
$ws.Range("E22").Value = "Harald"

# Investigation Person Email (row 24)
$ws.Range("B24").Value = "jens.krumsieck@thuenen.de"
$ws.Range("E24").Value = "harald.vonwaldow@thuenen.de"

# Investigation Person Affiliation (row 28)
$ws.Range("B28").Value = "Johann Heinrich von Thünen-Institut, Zentrum für Informationsmanagement;Technische Universität Braunschweig, Institut für Anorganische und Analytische Chemie"
$ws.Range("C28").Value = "Leibniz Institute of Plant Genetics and Crop Plant Research (IPK), Department of Breeding Research"
$ws.Range("D28").Value = "Leibniz Institute of Plant Genetics and Crop Plant Research (IPK), Department of Breeding Research"
$ws.Range("E28").Value = "Johann Heinrich von Thünen-Institut, Centre for Information Management;Eawag, IT Services"

# Investigation Person Roles (row 29)
$ws.Range("E29").Value = "Principal Investigator"

# Investigation Person Roles Term Accession Number (row 30)
$ws.Range("E30").Value = "NCIT:C19924"

# Investigation Person Roles Term Source REF (row 31)
$ws.Range("E31").Value = "NCIT"

# Comment[ORCID] (row 32)
$ws.Range("A32").Value = "Comment[ORCID]"
$ws.Range("B32").Value = "0000-0001-6242-5846"
$ws.Range("D32").Value = "0000-0002-8948-6793"
$ws.Range("E32").Value = "0000-0003-4800-2833"
